$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-64 down to 37-65
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new weekly data point
$ws.Range("A36").Value = 11
$ws.Range("B36").Value = "Vega Monumental Concepción"
$ws.Range("C36").Value = "Bíobío"
$ws.Range("D36").Value = 45167
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E36").Value = 8
$ws.Range("F36").Value = 100114007
$ws.Range("G36").Value = "Jengibre"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 20
$ws.Range("K36").Value = 14000
$ws.Range("L36").Value = 14000
$ws.Range("M36").Value = 14000
$ws.Range("N36").Value = "$/caja 13 kilos"
$ws.Range("O36").Value = "Perú"
$ws.Range("P36").Value = 1077
$ws.Range("Q36").Value = 13
$ws.Range("R36").Value = "Hortaliza"
